# CP_1.21mW_polarimeter.xlsx edit
#
# Commit: "Added min and max purity cases to computeRho.m and implemented a
#          scheme for plotting the coherency matrices onto the Bloch sphere
#          using Stokes parameters"
#
# Concretely this adds two new worksheets ("rho_min", "rho_max") after the
# existing "rho_mat" sheet, each holding the same theta/Jxx/Jyy/beta/gamma/
# trace_sq layout as "rho_mat" but computed from the min- and max-purity
# coherency-matrix cases, and it refreshes a handful of "rho_mat" values
# that shifted in the last significant digit once the min/max computation
# was folded into the same script.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Tiny last-digit refreshes on the existing "rho_mat" sheet.
# ---------------------------------------------------------------------
$rhoMat = $wb.Worksheets.Item("rho_mat")
$rhoMat.Cells.Item(3, 2).Value = 0.046448475234149959
$rhoMat.Cells.Item(4, 3).Value = 0.8170073770147358
$rhoMat.Cells.Item(6, 4).Value = -0.13156991358567202
$rhoMat.Cells.Item(6, 5).Value = -0.48138344766786306
$rhoMat.Cells.Item(11, 2).Value = 0.0046085791735393468
$rhoMat.Cells.Item(12, 2).Value = 0.057787635148832728
$rhoMat.Cells.Item(13, 4).Value = -0.29563455791987325
$rhoMat.Cells.Item(14, 4).Value = -0.27124591275110421

# ---------------------------------------------------------------------
# 2) New sheet data: min-purity case ("rho_min") and max-purity case
#    ("rho_max"). Each row is theta, Jxx, Jyy, beta, gamma, trace_sq.
# ---------------------------------------------------------------------
$headers = @("theta", "Jxx", "Jyy", "beta", "gamma", "trace_sq")

$rhoMinData = @(
    @(0, 0.10751885112282542, 0.89248114887717456, -0.010596991019683713, 0.047018582003302521, 1),
    @(10, 0.13237696944584326, 0.86762303055415679, -0.14285830967334529, -0.093219174465882682, 1),
    @(20, 0.24328480732079821, 0.75671519267920173, -0.22015329287319579, -0.22265936350511908, 1),
    @(30, 0.38205654323897253, 0.61794345676102747, -0.20432566090514426, -0.32547896278906718, 1),
    @(40, 0.47508718504258246, 0.52491281495741748, -0.10582651390917869, -0.3871943876219246, 1),
    @(50, 0.48375944976013191, 0.51624055023986815, 0.028141625455643415, -0.39626942585278141, 1),
    @(60, 0.39936054132861215, 0.60063945867138779, 0.1370308689327645, -0.36076150209992086, 1),
    @(70, 0.27635297324174463, 0.72364702675825532, 0.15694539481340974, -0.28472209584636626, 1),
    @(80, 0.15863259598023485, 0.84136740401976506, 0.09755378681666721, -0.17697967630100184, 0.99999999999999978),
    @(90, 0.10311053966080717, 0.89688946033919281, -0.029276244501280455, -0.045687462349127086, 1),
    @(100, 0.1412879993239885, 0.8587120006760115, -0.16376507238135107, 0.094911934360684666, 1),
    @(110, 0.26104696909219521, 0.73895303090780473, -0.2417091554116402, 0.22714153915719767, 1),
    @(120, 0.40122495728848212, 0.59877504271151794, -0.22250285156479652, 0.33008898235990136, 1),
    @(130, 0.4967685669901174, 0.5032314330098826, -0.11663903887009502, 0.38431314289981888, 1),
    @(140, 0.50325251847504326, 0.49674748152495674, 0.026739759587709609, 0.39626199744679336, 1),
    @(150, 0.42307029927008666, 0.57692970072991334, 0.13405121152965521, 0.36083349310465906, 1),
    @(160, 0.28288374549416007, 0.71711625450583993, 0.16762541330327385, 0.27910288245314552, 1),
    @(170, 0.15823002616890203, 0.84176997383109786, 0.10678235509483582, 0.16994591801548967, 0.99999999999999978),
    @(180, 0.10332930312388966, 0.89667069687611034, -0.010390213092367789, 0.044990134523109786, 1)
)

$rhoMaxData = @(
    @(0, 0.0037280465348734876, 0.99627195346512643, -0.013399342743558416, 0.059452545954612121, 0.99999999999999978),
    @(10, 0.046448475500789906, 0.9535515244992101, -0.17625012242039304, -0.11500829702596792, 1),
    @(20, 0.18299262328529478, 0.81700737671470514, -0.27185854144615335, -0.27495318835269816, 0.99999999999999978),
    @(30, 0.35330061352512782, 0.64669938647487224, -0.2541425350673801, -0.40483436269277939, 1),
    @(40, 0.46902688382941682, 0.53097311617058318, -0.13156991350995073, -0.48138344739083849, 1),
    @(50, 0.47957673438743614, 0.52042326561256391, 0.035389434653450638, -0.49832768094655794, 1),
    @(60, 0.37382737162118568, 0.62617262837881427, 0.17179687898296683, -0.45229006136116223, 1),
    @(70, 0.21662222497073449, 0.78337777502926553, 0.19886173953650388, -0.36076452789075286, 1),
    @(80, 0.069740237018533768, 0.93025976298146629, 0.12295687490787335, -0.22306533277968371, 1),
    @(90, 0.0046085796420971816, 0.99539142035790285, -0.036542165503310998, -0.057026399356484125, 1),
    @(100, 0.057787635435619433, 0.94221236456438062, -0.20188602487326024, 0.11700537155132996, 1),
    @(110, 0.20773645934662813, 0.7922635406533719, -0.29563455755536411, 0.27781690071616005, 1),
    @(120, 0.37958659224060654, 0.62041340775939346, -0.27124591250594798, 0.40240062812092575, 1),
    @(130, 0.495977164927634, 0.504022835072366, -0.14520481004516397, 0.47843429998406578, 1),
    @(140, 0.50409455060018538, 0.49590544939981474, 0.03366231414479845, 0.49884875733251999, 1),
    @(150, 0.40201070597430849, 0.59798929402569156, 0.17074788354105627, 0.45961207329123593, 1),
    @(160, 0.22258914104271463, 0.77741085895728534, 0.21417608733794929, 0.35661157905935958, 1),
    @(170, 0.068849860138126953, 0.93115013986187312, 0.13470822734319882, 0.21439041440645615, 1),
    @(180, 0.0033534720060074109, 0.99664652799399256, -0.013008934862294255, 0.056329328787963796, 1)
)

# Column widths (Excel character units) for the two new sheets.
$rhoMinColWidths = @(5.02, 11.88, 11.88, 13.59, 13.59, 7.74)
$rhoMaxColWidths = @(5.02, 13.88, 11.88, 13.59, 13.59, 7.74)

function Add-RhoSheet {
    param(
        [string]$sheetName,
        [object[]]$dataRows,
        [double[]]$colWidths
    )

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $sheetName

    # Header row, text-formatted like the "rho_mat" header row.
    for ($c = 0; $c -lt $headers.Length; $c++) {
        $cell = $ws.Cells.Item(1, $c + 1)
        $cell.NumberFormat = "@"
        $cell.Value = $headers[$c]
    }

    # Data rows.
    for ($r = 0; $r -lt $dataRows.Length; $r++) {
        $row = $dataRows[$r]
        for ($c = 0; $c -lt $row.Length; $c++) {
            $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
        }
    }

    for ($c = 0; $c -lt $colWidths.Length; $c++) {
        $ws.Columns.Item($c + 1).ColumnWidth = $colWidths[$c]
    }

    return $ws
}

Add-RhoSheet "rho_min" $rhoMinData $rhoMinColWidths | Out-Null
Add-RhoSheet "rho_max" $rhoMaxData $rhoMaxColWidths | Out-Null
